$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 414, shifting existing rows 414:450 down to 415:451
$ws.Rows.Item(414).Insert()

# Populate the newly inserted row 414 with its data
$ws.Cells.Item(414, 1).Value = 10
$ws.Cells.Item(414, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(414, 3).Value = "La Araucanía"
$ws.Cells.Item(414, 4).Value = 45021
$ws.Cells.Item(414, 5).Value = 9
$ws.Cells.Item(414, 6).Value = 100114013
$ws.Cells.Item(414, 7).Value = "Zanahoria"
$ws.Cells.Item(414, 8).Value = "Sin especificar"
$ws.Cells.Item(414, 9).Value = "Primera"
$ws.Cells.Item(414, 10).Value = 125
$ws.Cells.Item(414, 11).Value = 6000
$ws.Cells.Item(414, 12).Value = 6000
$ws.Cells.Item(414, 13).Value = 6000
$ws.Cells.Item(414, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(414, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(414, 16).Value = 240
$ws.Cells.Item(414, 17).Value = 25
$ws.Cells.Item(414, 18).Value = "Hortaliza"
